# Add new dialogue rows (14-20) and update row 13 border styling,
# matching the "Chansey" character sheet upload diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: becomes a "section separator" row (bottom border) ---
# Clone the border/font formatting already used by row 6 (same visual group)
# so the workbook reuses the existing cellXfs/border entries instead of
# minting new ones.
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Rows 16 and 17 are also "separator" rows needing borders ---
# Row 16 matches the row-6/row-13 bottom-border style; row 17 matches the
# row-11 top+bottom border style. Clone those first (formats only), then
# fill in all the values below.
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row number column (B) for the new rows ---
$ws.Cells.Item(14, 2).Value = 107
$ws.Cells.Item(15, 2).Value = 110
$ws.Cells.Item(16, 2).Value = 113
$ws.Cells.Item(17, 2).Value = 67
$ws.Cells.Item(18, 2).Value = 37
$ws.Cells.Item(19, 2).Value = 40
$ws.Cells.Item(20, 2).Value = 43

# --- English column (C), written first so new shared strings are appended
#     in the same order as the source edit (C block, then D block, then E) ---
$ws.Cells.Item(14, 3).Value = " Uh... Is it true? Team [CS:X]Charm[CR] is\nvisiting the guild? Really-really?"
$ws.Cells.Item(15, 3).Value = " Oh, really! What should I do?\nWhat should I do?"
$ws.Cells.Item(16, 3).Value = " ...[K]But, anyways…"

# --- Translated/Russian column (D) ---
$ws.Cells.Item(14, 4).Value = " Эм... Правда? Команда [CS:X]Шарм[CR]\nпосетила гильдию? Честно-честно?"
$ws.Cells.Item(15, 4).Value = " О, вот это да! Что же делать?\nЧто же делать?"
$ws.Cells.Item(16, 4).Value = " ...[K]Но, так или иначе..."

# --- Converted column (E) ---
$ws.Cells.Item(14, 5).Value = " Üí... Ðñàâäà? Ëïíàîäà [CS:X]Šàñí[CR]\nðïòåóéìà ãéìûäéý? Œåòóîï-œåòóîï?"
$ws.Cells.Item(15, 5).Value = " Ï, âïó üóï äà! Œóï çå äåìàóû?\nŒóï çå äåìàóû?"
$ws.Cells.Item(16, 5).Value = " ...[K]Îï, óàë éìé éîàœå..."

# Second block (rows 17-20): English, then Russian, then Converted
$ws.Cells.Item(17, 3).Value = " The two of you can do this!\nI know it! It just takes effort!"
$ws.Cells.Item(18, 3).Value = " Thank you!"
$ws.Cells.Item(19, 3).Value = " The way I can be standing here,\njust as always... I owe that to you!"
$ws.Cells.Item(20, 3).Value = " Honestly, thank you! ♪"

$ws.Cells.Item(17, 4).Value = " У вас всё получится! Я уверена\nв этом! Просто нужно поднажать!"
$ws.Cells.Item(18, 4).Value = " Спасибо вам!"
$ws.Cells.Item(19, 4).Value = " То, что я по прежнему могу быть\nздесь... За это я у вас в долгу!"
$ws.Cells.Item(20, 4).Value = " Честно, спасибо вам! ♪"

$ws.Cells.Item(17, 5).Value = " Ô âàò âòæ ðïìôœéóòÿ! Ÿ ôâåñåîà\nâ üóïí! Ðñïòóï îôçîï ðïäîàçàóû!"
$ws.Cells.Item(18, 5).Value = " Òðàòéáï âàí!"
$ws.Cells.Item(19, 5).Value = " Óï, œóï ÿ ðï ðñåçîåíô íïãô áúóû\nèäåòû... Èà üóï ÿ ô âàò â äïìãô!"
$ws.Cells.Item(20, 5).Value = " Œåòóîï, òðàòéáï âàí! ♪"

# --- Row heights (matches ht= on the corresponding rows in the sheet XML) ---
$ws.Rows.Item(14).RowHeight = 31.8
$ws.Rows.Item(15).RowHeight = 21.6
$ws.Rows.Item(17).RowHeight = 21.6
$ws.Rows.Item(19).RowHeight = 21.6

# --- View: scroll window down to the newly added rows, select D21 (the next empty data cell) ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D21").Select() | Out-Null

